# Update simulation for 5 generations (2018 -> 2019/2020 shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet from "2024" to "2019"
$ws.Name = "2019"

# 2. Update header text in K1 ("2018年家系号" -> "2019年家系号")
$ws.Range("K1").Value = "2019年家系号"

# 3. Update column H values: 2018XXX -> 2020XXX for rows 2..331
for ($row = 2; $row -le 331; $row++) {
    $seq = $row - 1
    $newVal = 2020000 + $seq
    $ws.Cells.Item($row, 8).Value = $newVal
}

# 4. Move the active selection from A1 to G9
$ws.Range("G9").Select()
